$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = ""

$ws.Range("H92").Value = 1124.9546
$ws.Range("I92").Value = 458.625
$ws.Range("K92").Value = 458.625
$ws.Range("M92").Value = 789.375

$ws.Range("H96").Value = 671
$ws.Range("I96").Value = 475.33334
$ws.Range("J96").Value = 866.6667
$ws.Range("K96").Value = 1426.00002
$ws.Range("L96").Value = 2600.0001
$ws.Range("M96").Value = -53.00001999999995
$ws.Range("N96").Value = -5346.0001

$ws.Range("H99").Value = 1327.0769
$ws.Range("I99").Value = 441.5
$ws.Range("J99").Value = 2086.1428
$ws.Range("K99").Value = 1324.5
$ws.Range("L99").Value = 6258.428400000001
$ws.Range("M99").Value = 173.5
$ws.Range("N99").Value = -9254.428400000001

$ws.Range("H100").Value = 2827.7778
$ws.Range("I100").Value = 2814
$ws.Range("J100").Value = 2845
$ws.Range("K100").Value = 2814
$ws.Range("L100").Value = 2845
$ws.Range("M100").Value = -2273
$ws.Range("N100").Value = -3927

$ws.Range("H112").Value = 23811034
$ws.Range("I112").Value = 200000270
$ws.Range("J112").Value = 1677.5946
$ws.Range("K112").Value = 600000810
$ws.Range("L112").Value = 5032.783799999999
$ws.Range("M112").Value = -599999702
$ws.Range("N112").Value = -7248.783799999999

$ws.Range("H113").Value = 2929.55
$ws.Range("J113").Value = 2816
$ws.Range("L113").Value = 2816
$ws.Range("N113").Value = -9324

$ws.Range("H132").Value = 6130449.5
$ws.Range("I132").Value = 7919
$ws.Range("J132").Value = 16334667
$ws.Range("K132").Value = 23757
$ws.Range("L132").Value = 49004001
$ws.Range("M132").Value = -21227
$ws.Range("N132").Value = -49009061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16051.267
$ws.Range("I32").Value = 17150.484
$ws.Range("K32").Value = 17150.484
$ws.Range("M32").Value = -16863.484

$ws.Range("H61").Value = 200401950
$ws.Range("I61").Value = 333668260
$ws.Range("J61").Value = 502500
$ws.Range("K61").Value = 333668260
$ws.Range("L61").Value = 502500
$ws.Range("M61").Value = -333668048
$ws.Range("N61").Value = -502924

$ws.Range("H74").Value = 5729369.5
$ws.Range("I74").Value = 7845489.5
$ws.Range("J74").Value = 86383.336
$ws.Range("K74").Value = 7845489.5
$ws.Range("L74").Value = 86383.336
$ws.Range("M74").Value = -7844615.5
$ws.Range("N74").Value = -88131.336

$ws.Range("H77").Value = 5729369.5
$ws.Range("I77").Value = 7845489.5
$ws.Range("J77").Value = 86383.336
$ws.Range("K77").Value = 39227447.5
$ws.Range("L77").Value = 431916.68
$ws.Range("M77").Value = -39223079.5
$ws.Range("N77").Value = -440652.68

$ws.Range("H132").Value = 125182
$ws.Range("I132").Value = 126312.625
$ws.Range("J132").Value = 124177
$ws.Range("K132").Value = 378937.875
$ws.Range("L132").Value = 372531
$ws.Range("M132").Value = -376407.875
$ws.Range("N132").Value = -377591

$ws.Range("H136").Value = 200401950
$ws.Range("I136").Value = 333668260
$ws.Range("J136").Value = 502500
$ws.Range("K136").Value = 1001004780
$ws.Range("L136").Value = 1507500
$ws.Range("M136").Value = -1001002230
$ws.Range("N136").Value = -1512600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1161.7222
$ws.Range("I20").Value = 950.1667
$ws.Range("J20").Value = 1584.8334
$ws.Range("K20").Value = 950.1667
$ws.Range("L20").Value = 1584.8334
$ws.Range("M20").Value = -703.1667
$ws.Range("N20").Value = -2078.8334

$ws.Range("H42").Value = 130000
$ws.Range("J42").Value = 130000
$ws.Range("L42").Value = 130000
$ws.Range("N42").Value = -130656

$ws.Range("H86").Value = 9080.031000000001
$ws.Range("I86").Value = 13876.944
$ws.Range("J86").Value = 2912.5715
$ws.Range("K86").Value = 13876.944
$ws.Range("L86").Value = 2912.5715
$ws.Range("M86").Value = -12753.944
$ws.Range("N86").Value = -5158.5715

$ws.Range("H89").Value = 9080.031000000001
$ws.Range("I89").Value = 13876.944
$ws.Range("J89").Value = 2912.5715
$ws.Range("K89").Value = 69384.72
$ws.Range("L89").Value = 14562.8575
$ws.Range("M89").Value = -63768.72
$ws.Range("N89").Value = -25794.8575

$ws.Range("H94").Value = 613.4828
$ws.Range("I94").Value = 625.9474
$ws.Range("J94").Value = 589.8
$ws.Range("K94").Value = 625.9474
$ws.Range("L94").Value = 589.8
$ws.Range("M94").Value = -174.9474
$ws.Range("N94").Value = -1491.8

$ws.Range("I105").Value = 35715970
$ws.Range("K105").Value = 35715970
$ws.Range("M105").Value = -35714223

$ws.Range("H134").Value = 4755
$ws.Range("I134").Value = 3981.111
$ws.Range("J134").Value = 5915.8335
$ws.Range("K134").Value = 11943.333
$ws.Range("L134").Value = 17747.5005
$ws.Range("M134").Value = -9408.332999999999
$ws.Range("N134").Value = -22817.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1700.4783
$ws.Range("I31").Value = 1221.6316
$ws.Range("J31").Value = 3975
$ws.Range("K31").Value = 1221.6316
$ws.Range("L31").Value = 3975
$ws.Range("M31").Value = -926.6315999999999
$ws.Range("N31").Value = -4565

$ws.Range("H34").Value = 1700.4783
$ws.Range("I34").Value = 1221.6316
$ws.Range("J34").Value = 3975
$ws.Range("K34").Value = 1221.6316
$ws.Range("L34").Value = 3975
$ws.Range("M34").Value = -1019.6316
$ws.Range("N34").Value = -4379

$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250

$ws.Range("H51").Value = 12074.25
$ws.Range("J51").Value = 12074.25
$ws.Range("L51").Value = 12074.25
$ws.Range("N51").Value = -13546.25

$ws.Range("H58").Value = 28573246
$ws.Range("I58").Value = 38462732
$ws.Range("J58").Value = 3624.3333
$ws.Range("K58").Value = 38462732
$ws.Range("L58").Value = 3624.3333
$ws.Range("M58").Value = -38462529
$ws.Range("N58").Value = -4030.3333

$ws.Range("H59").Value = 26194.154
$ws.Range("J59").Value = 26194.154
$ws.Range("L59").Value = 26194.154
$ws.Range("N59").Value = -28484.154

$ws.Range("H60").Value = 11256.728
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 11882.4
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 11882.4
$ws.Range("M60").Value = -4489
$ws.Range("N60").Value = -12904.4

$ws.Range("H61").Value = 12074.25
$ws.Range("J61").Value = 12074.25
$ws.Range("L61").Value = 12074.25
$ws.Range("N61").Value = -12770.25

$ws.Range("H105").Value = 777.9
$ws.Range("I105").Value = 753.2222
$ws.Range("K105").Value = 753.2222
$ws.Range("M105").Value = 993.7778

$ws.Range("H134").Value = 30237.9
$ws.Range("I134").Value = 2149.081
$ws.Range("K134").Value = 6447.243
$ws.Range("M134").Value = -3912.243

$ws.Range("H136").Value = 28573246
$ws.Range("I136").Value = 38462732
$ws.Range("J136").Value = 3624.3333
$ws.Range("K136").Value = 115388196
$ws.Range("L136").Value = 10872.9999
$ws.Range("M136").Value = -115385646
$ws.Range("N136").Value = -15972.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 185.44444
$ws.Range("I2").Value = 286.66666
$ws.Range("J2").Value = 134.83333
$ws.Range("K2").Value = 1719.99996
$ws.Range("L2").Value = 808.9999799999999
$ws.Range("M2").Value = -1606.99996
$ws.Range("N2").Value = -1034.99998

$ws.Range("H131").Value = 3194.0833
$ws.Range("I131").Value = 535.9
$ws.Range("J131").Value = 3893.6052
$ws.Range("K131").Value = 1607.7
$ws.Range("L131").Value = 11680.8156
$ws.Range("M131").Value = 3432.3
$ws.Range("N131").Value = -21760.8156

$ws.Range("H137").Value = 30316.715
$ws.Range("I137").Value = 1238
$ws.Range("J137").Value = 41948.2
$ws.Range("K137").Value = 3714
$ws.Range("L137").Value = 125844.6
$ws.Range("M137").Value = 1386
$ws.Range("N137").Value = -136044.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 46871.125
$ws.Range("I70").Value = 76336.28999999999
$ws.Range("J70").Value = 5619.9
$ws.Range("K70").Value = 76336.28999999999
$ws.Range("L70").Value = 5619.9
$ws.Range("M70").Value = -76066.28999999999
$ws.Range("N70").Value = -6159.9

$ws.Range("H73").Value = 46871.125
$ws.Range("I73").Value = 76336.28999999999
$ws.Range("J73").Value = 5619.9
$ws.Range("K73").Value = 76336.28999999999
$ws.Range("L73").Value = 5619.9
$ws.Range("M73").Value = -75400.28999999999
$ws.Range("N73").Value = -7491.9

$ws.Range("H122").Value = 1854
$ws.Range("I122").Value = 1821.5555
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5464.666499999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3014.666499999999
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 107646.52
$ws.Range("I132").Value = 101772.6
$ws.Range("J132").Value = 114173.11
$ws.Range("K132").Value = 305317.8
$ws.Range("L132").Value = 342519.33
$ws.Range("M132").Value = -302787.8
$ws.Range("N132").Value = -347579.33

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 36807.93
$ws.Range("I132").Value = 2105.28
$ws.Range("J132").Value = 253699.5
$ws.Range("K132").Value = 6315.84
$ws.Range("L132").Value = 761098.5
$ws.Range("M132").Value = -3785.84
$ws.Range("N132").Value = -766158.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 111558
$ws.Range("I100").Value = 71860.28999999999
$ws.Range("J100").Value = 250500
$ws.Range("K100").Value = 143720.58
$ws.Range("L100").Value = 501000
$ws.Range("M100").Value = -143179.58
$ws.Range("N100").Value = -502082

$ws.Range("H135").Value = 64700
$ws.Range("J135").Value = 64700
$ws.Range("L135").Value = 64700
$ws.Range("N135").Value = -74840
